$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - values change, label stays the same
$ws.Range("B3").Value = 0.03189485801624538
$ws.Range("C3").Value = 0.03308258631898218
$ws.Range("D3").Value = 0.03364621646255368

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.03037899043311675
$ws.Range("C4").Value = 0.03319805781288295
$ws.Range("D4").Value = 0.03319560161109698

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.03069757065543506
$ws.Range("C5").Value = 0.02867785842402152
$ws.Range("D5").Value = 0.02752305608722286
